# Apply updated stock-report figures (quantities/values recomputed, two
# rows (11 & 12) swapped/corrected, and subtotal + grand total rows
# refreshed to match).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 59408
$ws.Range("C11").Value = "SIG-3W Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D11").Value = 388.17
$ws.Range("E11").Value = 463.78
$ws.Range("F11").Value = 26
$ws.Range("G11").Value = 10092.42
$ws.Range("B12").Value = 47438
$ws.Range("C12").Value = "SIG-3w Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D12").Value = 401.81
$ws.Range("E12").Value = 480.05
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 803.62
$ws.Range("F13").Value = 37
$ws.Range("G13").Value = 13463.56
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 6781.1
$ws.Range("B19").Value = 53803.07
$ws.Range("F122").Value = 3
$ws.Range("G122").Value = 157.2
$ws.Range("B126").Value = 2855.25
$ws.Range("F143").Value = 137
$ws.Range("G143").Value = 7289.77
$ws.Range("F147").Value = 276
$ws.Range("G147").Value = 4642.32
$ws.Range("F152").Value = 37
$ws.Range("G152").Value = 916.86
$ws.Range("F156").Value = 56
$ws.Range("G156").Value = 3720.64
$ws.Range("B160").Value = 113302.53
$ws.Range("F167").Value = 4
$ws.Range("G167").Value = 3251.32
$ws.Range("F170").Value = 34
$ws.Range("G170").Value = 24522.84
$ws.Range("B188").Value = 169499.14
$ws.Range("F252").Value = 6
$ws.Range("G252").Value = 111
$ws.Range("B259").Value = 3874.09
$ws.Range("F305").Value = 25
$ws.Range("G305").Value = 2592.5
$ws.Range("F312").Value = 191
$ws.Range("G312").Value = 6360.3
$ws.Range("B331").Value = 209304.48
$ws.Range("F340").Value = 120
$ws.Range("G340").Value = 9602.4
$ws.Range("F344").Value = 71
$ws.Range("G344").Value = 2820.12
$ws.Range("F346").Value = 179
$ws.Range("G346").Value = 28801.1
$ws.Range("F350").Value = 227
$ws.Range("G350").Value = 16732.17
$ws.Range("F361").Value = 4
$ws.Range("G361").Value = 451.36
$ws.Range("F364").Value = 22
$ws.Range("G364").Value = 3156.56
$ws.Range("F379").Value = 321
$ws.Range("G379").Value = 7421.52
$ws.Range("F381").Value = 195
$ws.Range("G381").Value = 21414.9
$ws.Range("F385").Value = 155
$ws.Range("G385").Value = 8807.1
$ws.Range("F389").Value = 151
$ws.Range("G389").Value = 19155.86
$ws.Range("F391").Value = 78
$ws.Range("G391").Value = 4487.34
$ws.Range("F393").Value = 20
$ws.Range("G393").Value = 1379.2
$ws.Range("F401").Value = 785
$ws.Range("G401").Value = 46079.5
$ws.Range("F404").Value = 29
$ws.Range("G404").Value = 6282.27
$ws.Range("F406").Value = 124
$ws.Range("G406").Value = 6755.52
$ws.Range("F415").Value = 5
$ws.Range("G415").Value = 415.5
$ws.Range("F417").Value = 629
$ws.Range("G417").Value = 107766.57
$ws.Range("F418").Value = 208
$ws.Range("G418").Value = 31443.36
$ws.Range("F421").Value = 1
$ws.Range("G421").Value = 160.32
$ws.Range("F428").Value = 62
$ws.Range("G428").Value = 1259.84
$ws.Range("F429").Value = 418
$ws.Range("G429").Value = 24858.46
$ws.Range("F431").Value = 209
$ws.Range("G431").Value = 19359.67
$ws.Range("B435").Value = 631604.33
$ws.Range("F437").Value = 115
$ws.Range("G437").Value = 21111.7
$ws.Range("F444").Value = 28
$ws.Range("G444").Value = 5832.4
$ws.Range("F445").Value = 17
$ws.Range("G445").Value = 3280.32
$ws.Range("F448").Value = 6
$ws.Range("G448").Value = 2156.58
$ws.Range("B453").Value = 101029.35
$ws.Range("F479").Value = 76
$ws.Range("G479").Value = 1625.64
$ws.Range("F480").Value = 48
$ws.Range("G480").Value = 1020
$ws.Range("F482").Value = 176
$ws.Range("G482").Value = 1870.88
$ws.Range("F484").Value = 157
$ws.Range("G484").Value = 3819.81
$ws.Range("F493").Value = 77
$ws.Range("G493").Value = 7160.23
$ws.Range("B507").Value = 118253.72
$ws.Range("F522").Value = 292
$ws.Range("G522").Value = 48469.08
$ws.Range("B524").Value = 145694.4
$ws.Range("F703").Value = 15
$ws.Range("G703").Value = 1223.4
$ws.Range("B704").Value = 16199.68
$ws.Range("F843").Value = 237
$ws.Range("G843").Value = 19329.72
$ws.Range("F846").Value = 65
$ws.Range("G846").Value = 10043.8
$ws.Range("F847").Value = 181
$ws.Range("G847").Value = 14762.36
$ws.Range("F848").Value = 331
$ws.Range("G848").Value = 44056.1
$ws.Range("F852").Value = 92
$ws.Range("G852").Value = 1998.24
$ws.Range("F853").Value = 174
$ws.Range("G853").Value = 6486.72
$ws.Range("F863").Value = 358
$ws.Range("G863").Value = 51552
$ws.Range("F865").Value = 240
$ws.Range("G865").Value = 28970.4
$ws.Range("F866").Value = 3
$ws.Range("G866").Value = 362.13
$ws.Range("B867").Value = 440478.23
$ws.Range("F910").Value = 28
$ws.Range("G910").Value = 730.52
$ws.Range("F912").Value = 1608
$ws.Range("G912").Value = 262280.88
$ws.Range("B918").Value = 293342.11
$ws.Range("B930").Value = 5547433.43
$ws.Range("B931").Value = 5547433.43
